$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.007.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.70"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6356"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +5.87%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2920"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07308"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07650"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.828.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.964"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6614"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.049"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008574"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.017.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.082.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.076"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.002"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.471"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1371"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.503"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.086"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.203"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.994"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05300"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7390"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.827"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.149"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.642"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.286.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.86%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01777"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.372"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8927"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.980.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5142"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.45%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.78%  "

$ws.Range("B50").Value = "XinFinNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07515"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.24%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.722"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.33%  "

Write-Host "Crypto price updates applied"
